$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 new columns (N:P) for the "...1" UDP/GroupMask/Compatibility
#    block, right before the existing PrefaultTime2 column.
# ---------------------------------------------------------------------------
$ws.Columns("N:P").Insert()

$ws.Range("N1").Value = "UDPPortNumber1"
$ws.Range("N2").Value = "'1025"
$ws.Range("N3").Value = "'1025"

$ws.Range("O1").Value = "GroupMaskID1"
$ws.Range("O2").Value = "'13"
$ws.Range("O3").Value = "'13"

$ws.Range("P1").Value = "Compatibility1"
$ws.Range("P2").Value = "'None"
$ws.Range("P3").Value = "'None"

# ---------------------------------------------------------------------------
# 2. Insert 3 new columns (T:V) for the "...2" UDP/GroupMask/Compatibility
#    block, right before the existing TimeMaster column (which, after the
#    first insert above, now lives at column T).
# ---------------------------------------------------------------------------
$ws.Columns("T:V").Insert()

$ws.Range("T1").Value = "UDPPortNumber2"
$ws.Range("T2").Value = "'1025"
$ws.Range("T3").Value = "'1025"

$ws.Range("U1").Value = "GroupMaskID2"
$ws.Range("U2").Value = "'13"
$ws.Range("U3").Value = "'13"

$ws.Range("V1").Value = "Compatibility2"
$ws.Range("V2").Value = "'None"
$ws.Range("V3").Value = "'None"

# ---------------------------------------------------------------------------
# 3. Append 2 new columns (AB:AC) with header-only data.
# ---------------------------------------------------------------------------
$ws.Range("AB1").Value = "Delay"
$ws.Range("AC1").Value = "NumberofTimes"

# ---------------------------------------------------------------------------
# 4. Update the view state (top-left cell, active selection) to match the
#    saved file as closely as this runtime's COM surface allows.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R2").Select() | Out-Null
